$p = $ppt.ActivePresentation
$p.Slides.Item(3).Delete()
